$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.973.36'
$ws.Range('E2').Value = '  +5.44%  '
$ws.Range('D3').Value = '2.363.52'
$ws.Range('E3').Value = '  +3.23%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'548.88"
$ws.Range('E5').Value = '  +2.95%  '
$ws.Range('D6').Value = "'133.04"
$ws.Range('E6').Value = '  +2.11%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = "'0.590"
$ws.Range('E8').Value = '  +2.11%  '
$ws.Range('D9').Value = '2.361.04'
$ws.Range('E9').Value = '  +3.13%  '
$ws.Range('E10').Value = '  +2.34%  '
$ws.Range('D11').Value = "'5.50"
$ws.Range('E11').Value = '  +1.60%  '
$ws.Range('E12').Value = '  +1.30%  '
$ws.Range('E13').Value = '  +2.27%  '
$ws.Range('D14').Value = "'24.05"
$ws.Range('E14').Value = '  +2.56%  '
$ws.Range('D15').Value = '2.785.27'
$ws.Range('E15').Value = '  +3.18%  '
$ws.Range('D16').Value = '60.873.90'
$ws.Range('E16').Value = '  +5.34%  '
$ws.Range('E17').Value = '  +2.27%  '
$ws.Range('D18').Value = '2.321.40'
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('D19').Value = "'10.73"
$ws.Range('E19').Value = '  +2.27%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = "'4.19"
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = "'6.90"
$ws.Range('E21').Value = '  +8.71%  '
$ws.Range('D22').Value = "'316.52"
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('D23').Value = "'0.998"
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').Value = "'63.56"
$ws.Range('E24').Value = '  +2.04%  '
$ws.Range('D25').Value = "'0.173"
$ws.Range('E25').Value = '  +4.92%  '
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').Value = "'8.01"
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('E28').Value = '  +6.66%  '
$ws.Range('E29').Value = '  +3.21%  '
$ws.Range('D30').Value = "'172.41"
$ws.Range('E30').Value = '  +1.26%  '
$ws.Range('D31').Value = '0.0₃0738'
$ws.Range('E31').Value = '  +3.41%  '
$ws.Range('D32').Value = "'1.16"
$ws.Range('E32').Value = '  +10.70%  '
$ws.Range('E33').Value = '  +3.16%  '
$ws.Range('E34').Value = '  +16.87%  '
$ws.Range('E35').Value = '  +1.73%  '
$ws.Range('D36').Value = "'18.10"
$ws.Range('E36').Value = '  +2.35%  '
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('E39').Value = '  +8.06%  '
$ws.Range('D40').Value = "'316.88"
$ws.Range('E40').Value = '  +10.44%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = "'38.31"
$ws.Range('E41').Value = '  +0.55%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = "'1.54"
$ws.Range('E42').Value = '  +4.11%  '
$ws.Range('D43').Value = "'143.89"
$ws.Range('E43').Value = '  +2.78%  '
$ws.Range('D44').Value = "'3.48"
$ws.Range('E44').Value = '  +2.44%  '
$ws.Range('D45').Value = "'0.0956"
$ws.Range('E45').Value = '  +1.00%  '
$ws.Range('D46').Value = "'19.46"
$ws.Range('E46').Value = '  +7.96%  '
$ws.Range('E47').Value = '  +1.27%  '
$ws.Range('D48').Value = "'0.566"
$ws.Range('E48').Value = '  +2.26%  '
$ws.Range('E49').Value = '  +2.54%  '
$ws.Range('D50').Value = '0.0₆0217'
$ws.Range('E50').Value = '  +7.96%  '
$ws.Range('D51').Value = "'11.04"
$ws.Range('E51').Value = '  +1.05%  '
